$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "29.825.51"
$ws.Range("E2").Value = "  -0.41%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.895.37"
$ws.Range("E3").Value = "  +0.06%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.000"
$ws.Range("E4").Value = "  -0.09%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "0.7978"
$ws.Range("E5").Value = "  -3.67%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "242.89"
$ws.Range("E6").Value = "  +0.46%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.9999"
$ws.Range("E7").Value = "  -0.13%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3170"
$ws.Range("E8").Value = "  -1.78%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "25.53"
$ws.Range("E9").Value = "  -3.89%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.07039"
$ws.Range("E10").Value = "  +0.26%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.08059"
$ws.Range("E11").Value = "  +0.32%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "0.7712"
$ws.Range("E12").Value = "  +3.15%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.881.46"
$ws.Range("E13").Value = "  -0.71%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.334"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "92.38"
$ws.Range("E15").Value = "  +0.03%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "29.845.50"
$ws.Range("E16").Value = "  -0.37%  "
$ws.Range("B17").Value = "Uniswap"
$ws.Range("C17").Value = "https://coinranking.com/coin/_H5FVG9iW+uniswap-uni"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "6.006"
$ws.Range("E17").Value = "  +1.87%  "
$ws.Range("B18").Value = "Avalanche"
$ws.Range("C18").Value = "https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax"
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "13.89"
$ws.Range("E18").Value = "  -1.00%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "244.41"
$ws.Range("E19").Value = "  -0.28%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.000007709"
$ws.Range("E20").Value = "  -0.68%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "8.346"
$ws.Range("E21").Value = "  +20.86%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "0.9996"
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "2.141.37"
$ws.Range("E23").Value = "  -0.09%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "1.000"
$ws.Range("E24").Value = "  -0.04%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.1657"
$ws.Range("E25").Value = "  +3.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.338"
$ws.Range("E26").Value = "  +1.69%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "166.21"
$ws.Range("E27").Value = "  -0.39%  "
$ws.Range("E28").Value = "  -0.81%  "
$ws.Range("E29").Value = "  -0.67%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "1.406"
$ws.Range("E30").Value = "  +2.97%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.537"
$ws.Range("E31").Value = "  +1.09%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "4.437"
$ws.Range("E32").Value = "  +4.07%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "0.05714"
$ws.Range("E33").Value = "  +1.27%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.050"
$ws.Range("E34").Value = "  -0.63%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.260"
$ws.Range("E35").Value = "  -1.11%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.7391"
$ws.Range("E36").Value = "  +1.14%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.9995"
$ws.Range("E37").Value = "  -0.08%  "
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = "2.628"
$ws.Range("E38").Value = "  -3.45%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.01911"
$ws.Range("E39").Value = "  +0.05%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.787"
$ws.Range("E40").Value = "  +0.21%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.4412"
$ws.Range("E41").Value = "  -0.02%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "72.55"
$ws.Range("E42").Value = "  +0.80%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "5.822"
$ws.Range("E43").Value = "  -2.08%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.8453"
$ws.Range("E44").Value = "  +0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "0.9994"
$ws.Range("E45").Value = "  -0.11%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.035.69"
$ws.Range("E46").Value = "  +4.50%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "102.83"
$ws.Range("E47").Value = "  +2.12%  "
$ws.Range("E48").Value = "  -0.55%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "9.974"
$ws.Range("E49").Value = "  +3.12%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "7.441"
$ws.Range("E50").Value = "  -1.90%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "2.040.70"
$ws.Range("E51").Value = "  -0.12%  "
